$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The source-links table (Table1) currently spans A1:D33 - grow it by one
# row so the new entry becomes part of the table (ref + autofilter expand
# automatically).
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# New row of data: BBC article about plastic pollution & turtles.
$ws.Range("A34").Value = "https://www.bbc.co.uk/news/science-environment-51804884"
$ws.Range("B34").Value = "BBC Article: ""Why plastic is a deadly attraction for sea turtles"""
$ws.Range("C34").Value = "non-scientfic article. Motivations."
$ws.Range("D34").Value = "Roshi"

# Give the new row the same look as the row above it (link style in col A,
# wrapped plain text in B:D) before wiring up the hyperlink.
$ws.Range("A33:D33").Copy()
$ws.Range("A34:D34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Turn the URL text in A34 into a real hyperlink, then re-apply the link
# cell style (Hyperlinks.Add can nudge formatting) so it matches A2:A33.
$ws.Hyperlinks.Add($ws.Range("A34"), "https://www.bbc.co.uk/news/science-environment-51804884")
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the author's final selection/view state.
$ws.Range("C35").Select()
